# issue #5: stock data output to json file
#
# The 股票 (stock) sheet (sheet index 5 / "股票") gains a new
# "property_category" column holding the literal value "stock" for every
# data row. The column is inserted before the existing "date" column, so
# date / legislator_name / legislator_id all shift one column to the right
# (H->I, I->J, J->K).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at H, pushing date/legislator_name/legislator_id
# (old H/I/J) one column to the right (new I/J/K), carrying their values
# and styles along with them.
$ws.Columns.Item(8).Insert()

# New header for the inserted column.
$ws.Range("H1").Value = "property_category"

# New value for every existing data row (rows 2-4 hold the stock records).
$ws.Range("H2").Value = "stock"
$ws.Range("H3").Value = "stock"
$ws.Range("H4").Value = "stock"
